# Add "Save" column (H) to the s_vals sheet, matching the formatting of
# the existing header row (column G's style).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting from the G1 header cell onto H1, then set its text.
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)
$ws.Range("H1").Value = "Save"

# Fill in the Save values for rows 2-9 as per the source data.
$saveValues = @(0, 0, 1, 1, 0, 1, 1, 0)

for ($i = 0; $i -lt $saveValues.Length; $i++) {
    $row = 2 + $i
    $ws.Cells.Item($row, 8).Value = $saveValues[$i]
}
